$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 11, shifting the existing rows 11..42 down to 12..43
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly data point.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R carry the same constant
# values as the rest of the "Vega Monumental Concepción - Jengibre" block.
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Vega Monumental Concepción"
$ws.Range("C11").Value = "Bíobío"
$ws.Range("D11").Value = 44838
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 100114007
$ws.Range("G11").Value = "Jengibre"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 40
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("N11").Value = "$/caja 13 kilos"
$ws.Range("O11").Value = "Perú"
$ws.Range("P11").Value = 1115
$ws.Range("Q11").Value = 13
$ws.Range("R11").Value = "Hortaliza"
